$d = $word.ActiveDocument

# --- 1. Add alt text (descr) to the "Save to MP3" heading icon (docPr id=3 / InlineShape 10) ---
$d.InlineShapes.Item(10).AlternativeText = "Save To MP3 Button. (Keyboard shortcut: Ctrl M)."

# --- 2. Remove alt text (descr) from the save_to_mp3.png screenshot (docPr id=13 / InlineShape 11) ---
$d.InlineShapes.Item(11).AlternativeText = ""

# --- 3. Move the "_GoBack" bookmark from its old spot (inside "Keyboard") to just before the
#        final "M)." in "Ctrl Shift M). " -- this splits that run in two around the bookmark. ---
$rng = $d.Content
$rng.Find.Execute("Ctrl Shift M). ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.Start + "Ctrl Shift ".Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 4. Re-merge the "Open a Word Document..." runs that used to be split around the old
#        "_GoBack" bookmark location, now that the bookmark has moved away from here. ---
$openRng = $d.Content
$openRng.Find.Execute("Open a Word Document in the File menu (Keyboard shortcut: Ctrl O).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$openRng.Text = ""
$openRng.InsertAfter("Open a Word Document in the File menu (Keyboard shortcut: Ctrl O).")
